# Append a new data row (row 4) to the "Artfynd" export sheet, mirroring
# the existing rows 2-3 layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 131113864
$ws.Range("B4").Value = 99013
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("P4").Value = "Finnforsån, Finnforsån, Vb"
$ws.Range("Q4").Value = 754134
$ws.Range("R4").Value = 7191377
$ws.Range("S4").Value = 20
$ws.Range("T4").Value = "Västerbotten"
$ws.Range("U4").Value = "Skellefteå"
$ws.Range("V4").Value = "Västerbotten"
$ws.Range("W4").Value = "Skellefteå socken"

# Startdatum / Slutdatum are stored as plain text ("2025-08-14"), not real
# dates. Assigning that literal via .Value would get auto-recognised and
# converted into a date serial number, so instead stage the text (forced
# to Text via a leading apostrophe) in a scratch cell, copy/paste-special
# it into the two date columns (copy/paste preserves the Text type without
# touching cell formatting), then clear the scratch cell again.
$ws.Range("Z4").Value = "'2025-08-14"
$ws.Range("Z4").Copy()
$ws.Range("Y4").PasteSpecial()
$ws.Range("AA4").PasteSpecial()
$ws.Range("Z4").Clear()

$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "Emelie Bergkvist"
$ws.Range("AX4").Value = "Emelie Bergkvist"
